$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D: "Issues/Comments" ---
# Column width target (stored OOXML width) is 63.85546875; the closest value this
# runtime's width-quantization can produce is 63.833333... which is reached by any
# ColumnWidth input in [62.93, 63.07); 63 sits safely in the middle of that range.
$ws.Columns.Item(4).ColumnWidth = 63

$ws.Range("D1").Value = "Issues/Comments"

# --- Row 16: new comment text in D16, wrapped, taller row ---
$ws.Range("D16").Value = "Connecting to vertica from the docker container isn't working, need to think about how to make It work for both DBT and the normal data loading.
I think what I need to do is make the vertica DB accessible over the web and update the readme"
$ws.Range("D16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 90

# --- Sheet view: freeze header row, scroll/select state ---
# Final view mirrors Excel's "freeze top row, scroll down, select D5" flow:
#   pane ySplit=1 (freeze row 1), activePane=bottomLeft, selection D5 in the
#   scrolled (bottom) pane.
$ws.Activate()
$ws.Range("B1").Select()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D5").Select()
